$wb = $excel.ActiveWorkbook

# --- Sheets involved ---
$wsContactRRDUser = $wb.Worksheets.Item("Contact_RRDUser")   # sheet4.xml
$wsEditRRDTeam    = $wb.Worksheets.Item("Edit_RRDTeam")      # sheet6.xml
$wsEditRRDUser    = $wb.Worksheets.Item("Edit_RRDUser")      # sheet8.xml

# --- Value updates ---
# Edit_RRDTeam: G2 30 -> 10 (F2 stays 30)
$wsEditRRDTeam.Range("G2").Value = 10

# Edit_RRDUser: I2/J2/I3/J3 5 -> 15
$wsEditRRDUser.Range("I2").Value = 15
$wsEditRRDUser.Range("J2").Value = 15
$wsEditRRDUser.Range("I3").Value = 15
$wsEditRRDUser.Range("J3").Value = 15

# --- Apply the "quote prefix" number format (seen as a new cellXf with
# quotePrefix="1") to Edit_RRDTeam!F2:G2 and Edit_RRDUser!I2:J3, while
# keeping the cells' stored values numeric. Use a scratch cell entered
# with a leading apostrophe, then Paste Special (Formats only) onto the
# target ranges, so only the format/style is carried over.
$scratch = $wsEditRRDTeam.Range("Z1")
$scratch.Value = "'1"
$scratch.Copy()
$wsEditRRDTeam.Range("F2:G2").PasteSpecial(-4122)   # xlPasteFormats
$wsEditRRDUser.Range("I2:J3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$scratch.EntireColumn.Delete()

# --- Selection / active-cell changes ---
$wsContactRRDUser.Range("B3").Select()
$wsEditRRDTeam.Select()
$wsEditRRDTeam.Range("G2").Select()
